$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$desc = "We are seeking a Software Engineer to build and maintain high-quality software solutions.`nWork with global teams to drive innovation and deliver scalable applications.`nJoin Akkodis and be part of a tech-driven, collaborative environment."

$ws.Range("A16").Value = "JD_015"
$ws.Range("B16").Value = "Junior Dotnet Engineer"
$ws.Range("C16").Value = $desc
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 5
